$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet: handback status text + widened status columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet: fill in the handback target/file/datetime columns (I, J, K)
# ---------------------------------------------------------------------------
$wsZhCn.Range("J2").Value = "637ea68b-9530-496a-8634-572befa58fe0.b4dccdb3cbd5d4f2873307003f8f6c4628faa669.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-04 21:06:05"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19ea95bfe809148117ce50e52fbcadc0aa2a271a/e2e/637ea68b-9530-496a-8634-572befa58fe0.md", [Type]::Missing, [Type]::Missing, "637ea68b-9530-496a-8634-572befa58fe0.md")

$wsZhCn.Columns.Item(3).ColumnWidth = 29.14
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: fill in the handback target/file/datetime columns (I, J, K)
# ---------------------------------------------------------------------------
$wsDeDe.Range("J2").Value = "637ea68b-9530-496a-8634-572befa58fe0.b4dccdb3cbd5d4f2873307003f8f6c4628faa669.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-04 21:06:15"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19ea95bfe809148117ce50e52fbcadc0aa2a271a/e2e/637ea68b-9530-496a-8634-572befa58fe0.md", [Type]::Missing, [Type]::Missing, "637ea68b-9530-496a-8634-572befa58fe0.md")

$wsDeDe.Columns.Item(3).ColumnWidth = 29.14
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17

Write-Output "Handback report generated"
